# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 256
    $ws.Range("F5").Value = 6553
    $ws.Range("F6").Value = 5324
    $ws.Range("F7").Value = 443
}
